$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are plain text in the source data (e.g. "22.414.75",
# "1.002"), even when they look like numbers. Force text interpretation via
# NumberFormat "@" before assignment (otherwise COM auto-coerces numeric-looking
# strings to actual numbers), then restore the default "Normal" style so no
# stray number-format/style is left behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.414.75'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.572.68'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3768'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.83'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3422'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07655'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.015'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.936'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.573.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001134'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06769'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.221'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.402.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.422'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.733'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -8.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '146.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.038'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.32'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.743.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.193'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.008'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9907'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.03'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08594'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02547'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06572'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.333'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.470'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6451'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.55'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.306'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.086'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.38'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07334'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.71%  '
